$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New backlog entry appended after row 49 (written first so the new shared
# strings land in the same order as the source workbook: Admin Area /
# Update cars etc. before the two diary entries below).
$ws.Range("G51").Value = "Admin Area"
$ws.Range("H51").Value = "Update cars etc."

# Two new dated diary entries, written into the previously-blank row 35
# and the existing (until now H-only) row 36. Pull the date/description
# formatting from the row above (row 34) so the new cells pick up the
# same styles (date number format in B, Arial body font in C) instead of
# defaulting to General/Calibri.
$ws.Range("B34:C34").Copy()
$ws.Range("B35:C35").PasteSpecial(-4122)
$ws.Range("B36:C36").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B35").Value = 42785
$ws.Range("C35").Value = "Register working - user added to AspNetUser table"

$ws.Range("B36").Value = 42787
$ws.Range("C36").Value = "Checkout populated from View Model"

$ws.Range("C37").Select()
